# EspecificoCategoria_relaves.xlsx
# "Actualizacion de Utilidades y Contingencia Automatica"
#
# The "dedicacion" column (D) on Sheet1 was previously storing the text
# label "1.00" (a shared string) for every data row. It is updated here
# to hold the real numeric value 100 (e.g. 100%) instead, for all 55
# data rows (D2:D56).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = 100
}

# Reflect the author's final on-screen selection: column D (the column
# that was just edited) selected for the data rows.
$ws.Range("D2:D56").Select()
